# Appends newly-logged sensor rows (2026-01-28 afternoon readings) to the
# PIR, Humidity and Temperature sheets, matching the upstream CSV export.
# Each row tuple is: (rowNumber, ColA Date, ColB Timestamp, ColC Hour,
#                     ColD Location, ColE Value, ColF Status)

$wb = $excel.ActiveWorkbook

$pirRows = @(
    @("200","2026-01-28","16:27:10","16:00","Bathroom","No Motion","Inactive"),
    @("201","2026-01-28","16:27:11","16:00","Bathroom","No Motion","Inactive"),
    @("202","2026-01-28","16:27:13","16:00","Bathroom","No Motion","Inactive"),
    @("203","2026-01-28","16:27:18","16:00","Bathroom","No Motion","Inactive"),
    @("204","2026-01-28","16:27:24","16:00","Bathroom","No Motion","Inactive"),
    @("205","2026-01-28","16:27:28","16:00","Bathroom","No Motion","Inactive"),
    @("206","2026-01-28","16:27:33","16:00","Bathroom","No Motion","Inactive"),
    @("207","2026-01-28","16:27:38","16:00","Bathroom","No Motion","Inactive"),
    @("208","2026-01-28","16:27:44","16:00","Bathroom","No Motion","Inactive"),
    @("209","2026-01-28","16:27:48","16:00","Bathroom","No Motion","Inactive"),
    @("210","2026-01-28","16:27:54","16:00","Bathroom","No Motion","Inactive"),
    @("211","2026-01-28","16:27:58","16:00","Bathroom","No Motion","Inactive"),
    @("212","2026-01-28","16:28:04","16:00","Bathroom","No Motion","Inactive"),
    @("213","2026-01-28","16:28:09","16:00","Bathroom","No Motion","Inactive")
)

$humidityRows = @(
    @("199","2026-01-28","16:27:10","16:00","Bathroom","87.3%","Active"),
    @("200","2026-01-28","16:27:12","16:00","Bathroom","88.2%","Active"),
    @("201","2026-01-28","16:27:14","16:00","Bathroom","88.2%","Active"),
    @("202","2026-01-28","16:27:23","16:00","Bathroom","86.7%","Active"),
    @("203","2026-01-28","16:27:31","16:00","Bathroom","87.1%","Active"),
    @("204","2026-01-28","16:27:35","16:00","Bathroom","88.0%","Active"),
    @("205","2026-01-28","16:27:39","16:00","Bathroom","87.1%","Active"),
    @("206","2026-01-28","16:27:43","16:00","Bathroom","88.1%","Active"),
    @("207","2026-01-28","16:27:51","16:00","Bathroom","88.0%","Active"),
    @("208","2026-01-28","16:27:55","16:00","Bathroom","88.0%","Active"),
    @("209","2026-01-28","16:27:59","16:00","Bathroom","87.1%","Active"),
    @("210","2026-01-28","16:28:03","16:00","Bathroom","88.0%","Active")
)

$temperatureRows = @(
    @("198","2026-01-28","16:27:09","16:00","Bathroom","22.8C","Active"),
    @("199","2026-01-28","16:27:11","16:00","Bathroom","22.8C","Active"),
    @("200","2026-01-28","16:27:13","16:00","Bathroom","22.8C","Active"),
    @("201","2026-01-28","16:27:15","16:00","Bathroom","22.8C","Active"),
    @("202","2026-01-28","16:27:23","16:00","Bathroom","22.8C","Active"),
    @("203","2026-01-28","16:27:31","16:00","Bathroom","22.8C","Active"),
    @("204","2026-01-28","16:27:35","16:00","Bathroom","22.8C","Active"),
    @("205","2026-01-28","16:27:40","16:00","Bathroom","22.8C","Active"),
    @("206","2026-01-28","16:27:43","16:00","Bathroom","22.9C","Active"),
    @("207","2026-01-28","16:27:51","16:00","Bathroom","22.8C","Active"),
    @("208","2026-01-28","16:27:55","16:00","Bathroom","22.8C","Active"),
    @("209","2026-01-28","16:28:00","16:00","Bathroom","22.9C","Active"),
    @("210","2026-01-28","16:28:03","16:00","Bathroom","22.8C","Active")
)

function Add-LogRows {
    param($Worksheet, $Rows, $ValueColumnIsNumericLooking)

    foreach ($row in $Rows) {
        $r = [int]$row[0]

        # Column A holds a literal "yyyy-mm-dd" string in the source log, not
        # a real date -- force Text format first so Excel doesn't reinterpret
        # it as a date serial number.
        $Worksheet.Cells.Item($r, 1).NumberFormat = "@"
        $Worksheet.Cells.Item($r, 1).Value = $row[1]

        $Worksheet.Cells.Item($r, 2).Value = $row[2]
        $Worksheet.Cells.Item($r, 3).Value = $row[3]
        $Worksheet.Cells.Item($r, 4).Value = $row[4]

        if ($ValueColumnIsNumericLooking) {
            # Humidity's Value column is a "87.3%"-style string -- force Text
            # so it isn't converted into a numeric percentage.
            $Worksheet.Cells.Item($r, 5).NumberFormat = "@"
        }
        $Worksheet.Cells.Item($r, 5).Value = $row[5]

        $Worksheet.Cells.Item($r, 6).Value = $row[6]
    }
}

$pirSheet = $wb.Worksheets.Item("PIR")
Add-LogRows $pirSheet $pirRows $false

$humiditySheet = $wb.Worksheets.Item("Humidity")
Add-LogRows $humiditySheet $humidityRows $true

$temperatureSheet = $wb.Worksheets.Item("Temperature")
Add-LogRows $temperatureSheet $temperatureRows $false

